# Added additional CCS scenarios
# Fill in the new "connections" rows (14-22) describing the secondary/biofuel
# and fossil-fuel scenarios for the simplified steel factory model, and
# update workbook view state to match (connections tab active, new selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Data for the new / newly-populated rows. Columns: B,C,D,E,F,G,H,I
# (A, J, K, L stay empty on these rows, same as the existing rows 2-13.)
$rows = @(
    @{ r = 14; B = "steel";   C = "simple_BF";     D = "inflow"; E = "secondary fuel";    F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" },
    @{ r = 15; B = "steel";   C = "simple_BF";     D = "inflow"; E = "secondary biofuel"; F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" },
    @{ r = 16; B = "coke";    C = "simple_coke";   D = "inflow"; E = "fuel";              F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" },
    @{ r = 17; B = "power";   C = "simple_power";  D = "inflow"; E = "fuel";              F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" },
    @{ r = 18; B = "coke";    C = "simple_coke";   D = "inflow"; E = "biofuel";           F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" },
    @{ r = 19; B = "steel";   C = "simple_sinter"; D = "inflow"; E = "fossil fuel";       F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" },
    @{ r = 20; B = "steel";   C = "simple_sinter"; D = "inflow"; E = "biofuel";           F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" },
    @{ r = 21; B = "pellets"; C = "simple_pellets"; D = "inflow"; E = "biofuel";          F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" },
    @{ r = 22; B = "pellets"; C = "simple_pellets"; D = "inflow"; E = "fossil fuel";      F = "fuel"; G = "outflows"; H = "simple_fuel"; I = "fuel" }
)

# Rows 14-18 keep the "black font" style on column F that the blank
# placeholder rows already carried; rows 19-22 are plain (no explicit style).
$fStyledRows = @(14, 15, 16, 17, 18)

foreach ($row in $rows) {
    $r = $row.r

    foreach ($col in @("B", "C", "D", "E", "G")) {
        $cell = $ws.Range($col + $r)
        $cell.ClearFormats()
        $cell.Value = $row[$col]
    }

    $fCell = $ws.Range("F" + $r)
    $fCell.ClearFormats()
    $fCell.Value = $row["F"]
    if ($fStyledRows -contains $r) {
        $fCell.Font.Color = 0
    }

    $hCell = $ws.Range("H" + $r)
    $hCell.ClearFormats()
    $hCell.Value = $row["H"]
    $hCell.NumberFormat = "@"

    $iCell = $ws.Range("I" + $r)
    $iCell.ClearFormats()
    $iCell.Value = $row["I"]
    $iCell.Font.Color = 0
}

# The "connections" sheet becomes the active tab, with L28 selected.
$ws.Activate() | Out-Null
$ws.Range("L28").Select() | Out-Null
